$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Row 2 on each sheet keeps the same scenario text but loses the trailing
# "|PAYNOW" step from the ExecutionPipeline column (B).
$ws1.Range("B2").Value = "Search|AddToCart|CHECKOUTTRIP|LOGIN|ENTERPAXINFO|CONFIRMPAXINFO"
$ws2.Range("B2").Value = "Search|AddToCart|CHECKOUTTRIP|LOGIN|ENTERPAXINFO|CONFIRMPAXINFO"

# Only a single scenario (row 2) is kept per trip type; drop the other
# scenario rows (3-5) on both sheets so each table shrinks to A1:Q2.
$ws1.Rows("3:5").Delete()
$ws2.Rows("3:5").Delete()

# Leave sheet 2 ("Air_WorldSpan_RoundTrip") as the active tab/selection,
# and move sheet 1's selection off its old H3 cell.
$ws1.Range("B8").Select()
$ws2.Activate()
$ws2.Range("C6").Select()
